# Updated cryptos list values (prices & 1h volume %) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.896.97"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.912.19"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.52"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.80"
$ws.Range("E6").Value = "  +4.69%  "
$ws.Range("D7").Value = "3.912.01"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.39"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  +4.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.33"
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("D15").Value = "4.572.07"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "3.940.36"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "69.985.89"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.64"
$ws.Range("E18").Value = "  +8.95%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.16"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.88"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.746"
$ws.Range("E23").Value = "  +3.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000169"
$ws.Range("E24").Value = "  +6.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.73"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.33"
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.99"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.44"
$ws.Range("E31").Value = "  +2.82%  "
$ws.Range("D32").Value = "4.065.21"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.82"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.13"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "3.876.72"
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.13"
$ws.Range("E37").Value = "  +4.09%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.04"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.28"
$ws.Range("E40").Value = "  +10.50%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.330"
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.12"
$ws.Range("E43").Value = "  +7.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "436.35"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.19"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.66"
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0369"
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000274"
$ws.Range("E49").Value = "  +21.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.69"
$ws.Range("E50").Value = "  +5.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.50"
$ws.Range("E51").Value = "  -0.62%  "

# Restore default (no explicit) formatting on the cells we temporarily
# marked as Text so their saved style matches the source workbook.

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
